$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 62 (pushes old row 62 -> 63, old row 63 -> 64),
# carrying the existing row's formatting down with it.
$ws.Rows("62:62").Insert()

# Populate the newly inserted timesheet entry row (row 62).
$ws.Range("A62").Value = "Added Picasso page"
$ws.Range("B62").Value = 43758
$ws.Range("C62").Value = 1.5
$ws.Range("D62").Value = 25
$ws.Range("E62").Formula = "=D62*C62"

# Fix up the totals row (now row 64) so the SUM ranges include the new row.
$ws.Range("C64").Formula = "=SUM(C53:C62)"
$ws.Range("E64").Formula = "=SUM(E53:E62)"
